$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-extracted responses (replacing the old manual name/roll input) appended
# below the existing header + single response row.
$data = @(
    @("Holesh",        "12312",      "OBC", "Female", "Chhattisgarh",      298, 298, $null, $null, $null, $null),
    @("Holesh",        "12312",      "OBC", "Female", "Chhattisgarh",      298, 298, $null, $null, $null, $null),
    @("Holesh",        "1232141",    "SC",  "Male",   "Assam",             341, 75,  90,    123,   53,    26),
    @("HOLESH SHARMA", "3010018033", "UR",  "Male",   "Arunachal Pradesh", 341, 75,  90,    123,   53,    26)
)

$rowIndex = 3
foreach ($row in $data) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $value = $row[$col - 1]
        if ($null -ne $value) {
            $cell = $ws.Cells.Item($rowIndex, $col)
            # Roll numbers must stay text (leading content is numeric-looking
            # but the source sheet stores Roll as a string), so force the
            # cell to text format before writing it.
            if ($col -eq 2) {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $value
        }
    }
    $rowIndex++
}
